# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple sheets to match the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 176.43333
$ws.Range("I28").Value = 179
$ws.Range("K28").Value = 179
$ws.Range("M28").Value = 306
$ws.Range("H31").Value = 1495
$ws.Range("I31").Value = 1495
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 4485
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4255
$ws.Range("H53").Value = 17241918
$ws.Range("I53").Value = 38461708
$ws.Range("J53").Value = 838.5
$ws.Range("K53").Value = 38461708
$ws.Range("L53").Value = 838.5
$ws.Range("M53").Value = -38461071
$ws.Range("N53").Value = -2112.5
$ws.Range("H99").Value = 261.02777
$ws.Range("I99").Value = 238.65625
$ws.Range("J99").Value = 440
$ws.Range("K99").Value = 715.96875
$ws.Range("L99").Value = 1320
$ws.Range("M99").Value = 782.03125
$ws.Range("N99").Value = -4316
$ws.Range("H100").Value = 6439.3184
$ws.Range("I100").Value = 8974.23
$ws.Range("J100").Value = 2777.7778
$ws.Range("K100").Value = 8974.23
$ws.Range("L100").Value = 2777.7778
$ws.Range("M100").Value = -8433.23
$ws.Range("N100").Value = -3859.7778
$ws.Range("H129").Value = 1556
$ws.Range("J129").Value = 1361.75
$ws.Range("L129").Value = 4085.25
$ws.Range("N129").Value = -14085.25
$ws.Range("H135").Value = 764.24286
$ws.Range("I135").Value = 599.371
$ws.Range("J135").Value = 2042
$ws.Range("K135").Value = 5394.339
$ws.Range("L135").Value = 18378
$ws.Range("M135").Value = -2859.339
$ws.Range("N135").Value = -23448
$ws.Range("H141").Value = 3202.1892
$ws.Range("I141").Value = 1676.6428
$ws.Range("J141").Value = 7948.3335
$ws.Range("K141").Value = 5029.928400000001
$ws.Range("L141").Value = 23845.0005
$ws.Range("M141").Value = 150.0715999999993
$ws.Range("N141").Value = -34205.00049999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 802.6
$ws.Range("I32").Value = 790.5914
$ws.Range("J32").Value = 962.1429000000001
$ws.Range("K32").Value = 790.5914
$ws.Range("L32").Value = 962.1429000000001
$ws.Range("M32").Value = -503.5914
$ws.Range("N32").Value = -1536.1429
$ws.Range("H74").Value = 918.4286
$ws.Range("I74").Value = 918.4286
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 918.4286
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -44.42859999999996
$ws.Range("H77").Value = 918.4286
$ws.Range("I77").Value = 918.4286
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4592.143
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -224.143
$ws.Range("H110").Value = 1271.1333
$ws.Range("I110").Value = 854.3333
$ws.Range("J110").Value = 2938.3333
$ws.Range("K110").Value = 854.3333
$ws.Range("L110").Value = 2938.3333
$ws.Range("M110").Value = 1190.6667
$ws.Range("N110").Value = -7028.3333
$ws.Range("H132").Value = 1510205.2
$ws.Range("I132").Value = 1560
$ws.Range("K132").Value = 4680
$ws.Range("M132").Value = -2150

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1455452
$ws.Range("I86").Value = 2954.1667
$ws.Range("J86").Value = 2326950.8
$ws.Range("K86").Value = 2954.1667
$ws.Range("L86").Value = 2326950.8
$ws.Range("M86").Value = -1831.1667
$ws.Range("N86").Value = -2329196.8
$ws.Range("H89").Value = 1455452
$ws.Range("I89").Value = 2954.1667
$ws.Range("J89").Value = 2326950.8
$ws.Range("K89").Value = 14770.8335
$ws.Range("L89").Value = 11634754
$ws.Range("M89").Value = -9154.833500000001
$ws.Range("N89").Value = -11645986

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6130.1113
$ws.Range("I94").Value = 12905
$ws.Range("J94").Value = 710.2
$ws.Range("K94").Value = 12905
$ws.Range("L94").Value = 710.2
$ws.Range("M94").Value = -12454
$ws.Range("N94").Value = -1612.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 23421.553
$ws.Range("J12").Value = 31450.828
$ws.Range("L12").Value = 94352.484
$ws.Range("N12").Value = -94698.484

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 29419070
$ws.Range("I122").Value = 38470400
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 115411200
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -115408750
$ws.Range("N122").Value = -11650
$ws.Range("H132").Value = 9446.814
$ws.Range("I132").Value = 5857.6
$ws.Range("J132").Value = 19701.715
$ws.Range("K132").Value = 17572.8
$ws.Range("L132").Value = 59105.145
$ws.Range("M132").Value = -15042.8
$ws.Range("N132").Value = -64165.145
$ws.Range("H141").Value = 28750
$ws.Range("J141").Value = 38000
$ws.Range("L141").Value = 38000
$ws.Range("N141").Value = -48360

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1086.0769
$ws.Range("I22").Value = 398.42856
$ws.Range("J22").Value = 1339.421
$ws.Range("K22").Value = 398.42856
$ws.Range("L22").Value = 1339.421
$ws.Range("M22").Value = -103.42856
$ws.Range("N22").Value = -1929.421
$ws.Range("H27").Value = 1086.0769
$ws.Range("I27").Value = 398.42856
$ws.Range("J27").Value = 1339.421
$ws.Range("K27").Value = 398.42856
$ws.Range("L27").Value = 1339.421
$ws.Range("M27").Value = -291.42856
$ws.Range("N27").Value = -1553.421
$ws.Range("H40").Value = 6758605
$ws.Range("I40").Value = 1389.069
$ws.Range("K40").Value = 1389.069
$ws.Range("M40").Value = -1253.069
$ws.Range("H93").Value = 974.55554
$ws.Range("I93").Value = 810.6923
$ws.Range("J93").Value = 1400.6
$ws.Range("K93").Value = 810.6923
$ws.Range("L93").Value = 1400.6
$ws.Range("M93").Value = 437.3077
$ws.Range("N93").Value = -3896.6
$ws.Range("H136").Value = 34331604
$ws.Range("I136").Value = 11615951
$ws.Range("J136").Value = 500002500
$ws.Range("K136").Value = 34847853
$ws.Range("L136").Value = 1500007500
$ws.Range("M136").Value = -34845303
$ws.Range("N136").Value = -1500012600

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 12687.375
$ws.Range("J14").Value = 12687.375
$ws.Range("L14").Value = 12687.375
$ws.Range("N14").Value = -13023.375
$ws.Range("H126").Value = 2256.7144
$ws.Range("I126").Value = 1559.4
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 4678.200000000001
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -2208.200000000001
$ws.Range("N126").Value = -16940
$ws.Range("H136").Value = 10001241
$ws.Range("I136").Value = 13889618
$ws.Range("J136").Value = 2558.5715
$ws.Range("K136").Value = 41668854
$ws.Range("L136").Value = 7675.7145
$ws.Range("M136").Value = -41666304
$ws.Range("N136").Value = -12775.7145
$ws.Range("H140").Value = 38895.855
$ws.Range("J140").Value = 38895.855
$ws.Range("L140").Value = 38895.855
$ws.Range("N140").Value = -49255.855
$ws.Range("H141").Value = 45600
$ws.Range("I141").Value = 19500
$ws.Range("J141").Value = 58650
$ws.Range("K141").Value = 19500
$ws.Range("L141").Value = 58650
$ws.Range("M141").Value = -14320
$ws.Range("N141").Value = -69010

